# Update contribution percentages in the first "Names / Percentage
# Contribution" table:
#   - Jacob Artis:   9    -> 10
#   - Muqdas Sheikh: 13.6 -> 12.6 (written back out as three runs,
#                                   matching the author's edit history)

$d = $word.ActiveDocument
$table = $d.Tables.Item(1)

function Find-RowByName($tbl, $name) {
    for ($i = 1; $i -le $tbl.Rows.Count; $i++) {
        # Cell Range.Text carries trailing cell-mark characters (CR + BEL)
        # that Word appends to every cell's range text.
        $txt = $tbl.Cell($i, 1).Range.Text.TrimEnd([char]13, [char]7)
        if ($txt -eq $name) {
            return $i
        }
    }
    return -1
}

# --- Jacob Artis: 9 -> 10 -----------------------------------------------
$jacobRow = Find-RowByName $table "Jacob Artis"
$table.Cell($jacobRow, 2).Range.Text = "10"

# --- Muqdas Sheikh: 13.6 -> 12.6 -----------------------------------------
# Re-fetch the table/cell fresh (the previous edit shifted character
# offsets later in the document).
$table = $d.Tables.Item(1)
$muqdasRow = Find-RowByName $table "Muqdas Sheikh"
$cell = $table.Cell($muqdasRow, 2)
$cellStart = $cell.Range.Start

# Replace the cell text, then split it into three runs: "1", "2", ".6"
$cell.Range.Text = "12.6"

$seg1 = $d.Range($cellStart, $cellStart + 1)      # "1"
$seg2 = $d.Range($cellStart + 1, $cellStart + 2)  # "2"
$seg3 = $d.Range($cellStart + 2, $cellStart + 4)  # ".6"

# Give each segment distinct formatting first, forcing the engine to
# keep them as separate runs, then converge them back to the original
# 10pt size shared by the rest of the cell's text.
$seg1.Font.Size = 11
$seg2.Font.Size = 12
$seg3.Font.Size = 13

$seg1b = $d.Range($cellStart, $cellStart + 1)
$seg2b = $d.Range($cellStart + 1, $cellStart + 2)
$seg3b = $d.Range($cellStart + 2, $cellStart + 4)

$seg1b.Font.Size = 10
$seg2b.Font.Size = 10
$seg3b.Font.Size = 10
